$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update the command-line example text (drop --newSheet, split in/in -> old/new)
$ws.Range("B1").Value = "xltablediff.py --key=ID test1old.xlsx test1new.xlsx --out test1diff.xlsx"

# Row 4 used to hold the "Diff test: ..." text; it becomes a blank "+" marker row
# matching the style/content pattern of row 3 ("Example TableNew" / style "+").
$plusColor = $ws.Range("A3").Interior.Color
$ws.Range("A4:G4").Interior.Color = $plusColor
$ws.Range("A4").Value = "+"
$ws.Range("B4:G4").Value = ""

# Rows 5-7: update/add the test description text
$ws.Range("B5").Value = "Diff test:  xltablediff.py  --key=ID test1old.xlsx test1new.xlsx --out test1diff.xlsx"
$ws.Range("B6").Value = "Merge test:  xltablediff.py  --key=ID test1old.xlsx test1new.xlsx --merge=Color --out test1merge.xlsx"
$ws.Range("B7").Value = "Append test:  xltablediff.py  --key=ID test1old.xlsx test1new.xlsx --append --out test1append.xlsx"

# Insert a new blank row above the old header row (row 8), pushing the table
# (header + data rows) down by one; this also bumps the trailing rows down,
# which already lines up with the desired final content.
$ws.Rows.Item(8).Insert()

# The newly inserted row 8 becomes a blank "=" separator row.
# (leading apostrophe forces literal text so "=" isn't parsed as a formula)
$ws.Range("A8").Value = "'="
